$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9.084673430066717
$ws.Range("D2").Value = 4.848446813053727
$ws.Range("E2").Value = 11.51079637086398
$ws.Range("F2").Value = 28.30271294507854
$ws.Range("G2").Value = 3.6303306947679
$ws.Range("I2").Value = 25.24065257700906
$ws.Range("L2").Value = 8.980559138224811
$ws.Range("M2").Value = 23.53127622920516
$ws.Range("N2").Value = 17.38162619306493
$ws.Range("O2").Value = 24.88150095661491
$ws.Range("C3").Value = 9.128212910553495
$ws.Range("D3").Value = 4.863358258134007
$ws.Range("E3").Value = 11.58186531643598
$ws.Range("F3").Value = 28.04551728471385
$ws.Range("G3").Value = 3.633504299366391
$ws.Range("I3").Value = 25.15898211692503
$ws.Range("L3").Value = 9.023396058707156
$ws.Range("M3").Value = 22.82997654462756
$ws.Range("N3").Value = 17.10751983567498
$ws.Range("O3").Value = 24.75274205520819
$ws.Range("C4").Value = 9.15655344713956
$ws.Range("D4").Value = 4.872884476518452
$ws.Range("E4").Value = 11.627593229185
$ws.Range("F4").Value = 27.89593279868928
$ws.Range("G4").Value = 3.635556156253677
$ws.Range("I4").Value = 25.11637246226427
$ws.Range("L4").Value = 9.050939716598547
$ws.Range("M4").Value = 22.38875882112294
$ws.Range("N4").Value = 16.93897827465235
$ws.Range("O4").Value = 24.68104866438964
$ws.Range("C5").Value = 9.168506615005168
$ws.Range("D5").Value = 4.876860088995673
$ws.Range("E5").Value = 11.64675504714759
$ws.Range("F5").Value = 27.83713478905526
$ws.Range("G5").Value = 3.636418359645303
$ws.Range("I5").Value = 25.10091321589388
$ws.Range("L5").Value = 9.062477188343106
$ws.Range("M5").Value = 22.20653053259134
$ws.Range("N5").Value = 16.87031843293093
$ws.Range("O5").Value = 24.65370662048793
$ws.Range("C6").Value = 9.170515842009566
$ws.Range("D6").Value = 4.877525902232974
$ws.Range("E6").Value = 11.64996874403341
$ws.Range("F6").Value = 27.82750356138311
$ws.Range("G6").Value = 3.636563104219912
$ws.Range("I6").Value = 25.09846147532916
$ws.Range("L6").Value = 9.06441192285372
$ws.Range("M6").Value = 22.17613240494601
$ws.Range("N6").Value = 16.8589212886262
$ws.Range("O6").Value = 24.64928022379915
$ws.Range("C7").Value = 9.156713015205508
$ws.Range("D7").Value = 4.872937713483981
$ws.Range("E7").Value = 11.62784951515557
$ws.Range("F7").Value = 27.89513100739396
$ws.Range("G7").Value = 3.635567678613226
$ws.Range("I7").Value = 25.11615625157468
$ws.Range("L7").Value = 9.051094045514303
$ws.Range("M7").Value = 22.38631072715762
$ws.Range("N7").Value = 16.93805210445131
$ws.Range("O7").Value = 24.68067230786033
$ws.Range("C8").Value = 9.099352310302693
$ws.Range("D8").Value = 4.853511634911215
$ws.Range("E8").Value = 11.53486777227254
$ws.Range("F8").Value = 28.21234132545547
$ws.Range("G8").Value = 3.631403579262564
$ws.Range("I8").Value = 25.21093487580545
$ws.Range("L8").Value = 8.995072212149152
$ws.Range("M8").Value = 23.29181054609866
$ws.Range("N8").Value = 17.28721297163021
$ws.Range("O8").Value = 24.83559016500117
$ws.Range("C9").Value = 8.999617696864897
$ws.Range("D9").Value = 4.818337261983947
$ws.Range("E9").Value = 11.36905700313125
$ws.Range("F9").Value = 28.89752327721128
$ws.Range("G9").Value = 3.624052859998384
$ws.Range("I9").Value = 25.45606075038796
$ws.Range("L9").Value = 8.895019166101859
$ws.Range("M9").Value = 24.97361131821912
$ws.Range("N9").Value = 17.96647438990222
$ws.Range("O9").Value = 25.19672671084662
$ws.Range("C10").Value = 8.93411195055632
$ws.Range("D10").Value = 4.794246480330957
$ws.Range("E10").Value = 11.25721473150875
$ws.Range("F10").Value = 29.43511762592122
$ws.Range("G10").Value = 3.619143327436044
$ws.Range("I10").Value = 25.67140664495036
$ws.Range("L10").Value = 8.827422409977096
$ws.Range("M10").Value = 26.14041558990472
$ws.Range("N10").Value = 18.45760202179382
$ws.Range("O10").Value = 25.49534310028722
$ws.Range("C11").Value = 8.905998537805065
$ws.Range("D11").Value = 4.783661241954138
$ws.Range("E11").Value = 11.20848116797336
$ws.Range("F11").Value = 29.68612854299129
$ws.Range("G11").Value = 3.617015232887456
$ws.Range("I11").Value = 25.77679025085013
$ws.Range("L11").Value = 8.79794062486928
$ws.Range("M11").Value = 26.65420176273294
$ws.Range("N11").Value = 18.67838153040469
$ws.Range("O11").Value = 25.63800795968883
$ws.Range("C12").Value = 8.895595193160275
$ws.Range("D12").Value = 4.779706172462843
$ws.Range("E12").Value = 11.19033384507947
$ws.Range("F12").Value = 29.78202314867505
$ws.Range("G12").Value = 3.616224422143013
$ws.Range("I12").Value = 25.81773993398554
$ws.Range("L12").Value = 8.786957970798802
$ws.Range("M12").Value = 26.84616462390811
$ws.Range("N12").Value = 18.76153259206387
$ws.Range("O12").Value = 25.69297426175001
$ws.Range("C13").Value = 8.897824946218359
$ws.Range("D13").Value = 4.780555602206443
$ws.Range("E13").Value = 11.19422856137279
$ws.Range("F13").Value = 29.76133445701813
$ws.Range("G13").Value = 3.616394069289851
$ws.Range("I13").Value = 25.8088746955152
$ws.Range("L13").Value = 8.789315226129014
$ws.Range("M13").Value = 26.80493973442766
$ws.Range("N13").Value = 18.74364582987728
$ws.Range("O13").Value = 25.68109502075336
$ws.Range("C14").Value = 8.905137785031297
$ws.Range("D14").Value = 4.783334789579606
$ws.Range("E14").Value = 11.20698203100749
$ws.Range("F14").Value = 29.69400146873337
$ws.Range("G14").Value = 3.616949871229914
$ws.Range("I14").Value = 25.78013844382697
$ws.Range("L14").Value = 8.797033443899098
$ws.Range("M14").Value = 26.67004752175029
$ws.Range("N14").Value = 18.68523191653113
$ws.Range("O14").Value = 25.64251141708532
$ws.Range("C15").Value = 8.90964871010622
$ws.Range("D15").Value = 4.78504405488974
$ws.Range("E15").Value = 11.21483384438213
$ws.Range("F15").Value = 29.65286517103117
$ws.Range("G15").Value = 3.617292273694697
$ws.Range("I15").Value = 25.7626717658249
$ws.Range("L15").Value = 8.801784678979942
$ws.Range("M15").Value = 26.58707972126619
$ws.Range("N15").Value = 18.64939052667957
$ws.Range("O15").Value = 25.61899936065226
$ws.Range("C16").Value = 8.935983162901291
$ws.Range("D16").Value = 4.794945736165251
$ws.Range("E16").Value = 11.26044261748583
$ws.Range("F16").Value = 29.41883605379702
$ws.Range("G16").Value = 3.619284514474404
$ws.Range("I16").Value = 25.6646669163928
$ws.Range("L16").Value = 8.829374556641092
$ws.Range("M16").Value = 26.10648299088983
$ws.Range("N16").Value = 18.44311408259867
$ws.Range("O16").Value = 25.48615363136527
$ws.Range("C17").Value = 8.95257033915802
$ws.Range("D17").Value = 4.801115532836915
$ws.Range("E17").Value = 11.28897034097036
$ws.Range("F17").Value = 29.27685954176757
$ws.Range("G17").Value = 3.620533592028406
$ws.Range("I17").Value = 25.60642857376978
$ws.Range("L17").Value = 8.846624240410874
$ws.Range("M17").Value = 25.80718494503652
$ws.Range("N17").Value = 18.31584094486455
$ws.Range("O17").Value = 25.40637756135044
$ws.Range("C18").Value = 8.96226947436598
$ws.Range("D18").Value = 4.804699443551535
$ws.Range("E18").Value = 11.30558061720516
$ws.Range("F18").Value = 29.19581139567238
$ws.Range("G18").Value = 3.621261942638517
$ws.Range("I18").Value = 25.57363156231819
$ws.Range("L18").Value = 8.856665237422202
$ws.Range("M18").Value = 25.63344463259532
$ws.Range("N18").Value = 18.24239282978471
$ws.Range("O18").Value = 25.36113713742486
$ws.Range("C19").Value = 8.9655806733848
$ws.Range("D19").Value = 4.805918953695191
$ws.Range("E19").Value = 11.31123928469243
$ws.Range("F19").Value = 29.1684777722669
$ws.Range("G19").Value = 3.621510255052073
$ws.Range("I19").Value = 25.56264800482431
$ws.Range("L19").Value = 8.860085487807952
$ws.Range("M19").Value = 25.57435063556224
$ws.Range("N19").Value = 18.21748497301982
$ws.Range("O19").Value = 25.34593137857252
$ws.Range("C20").Value = 8.950788186808897
$ws.Range("D20").Value = 4.800455106739647
$ws.Range("E20").Value = 11.28591263132513
$ws.Range("F20").Value = 29.2919103625987
$ws.Range("G20").Value = 3.620399600133216
$ws.Range("I20").Value = 25.61255585012195
$ws.Range("L20").Value = 8.844775627801823
$ws.Range("M20").Value = 25.83921165415883
$ws.Range("N20").Value = 18.32941518600047
$ws.Range("O20").Value = 25.41480341544047
$ws.Range("C21").Value = 8.902983241122026
$ws.Range("D21").Value = 4.78251703083548
$ws.Range("E21").Value = 11.20322770561343
$ws.Range("F21").Value = 29.71375663507915
$ws.Range("G21").Value = 3.616786210891248
$ws.Range("I21").Value = 25.78855086738859
$ws.Range("L21").Value = 8.794761498871368
$ws.Range("M21").Value = 26.70974024480306
$ws.Range("N21").Value = 18.70240237643776
$ws.Range("O21").Value = 25.65381910539384
$ws.Range("C22").Value = 8.873154097856258
$ws.Range("D22").Value = 4.771104118602356
$ws.Range("E22").Value = 11.15097728498639
$ws.Range("F22").Value = 29.99432595585874
$ws.Range("G22").Value = 3.614512349855745
$ws.Range("I22").Value = 25.90964173994733
$ws.Range("L22").Value = 8.763131625288878
$ws.Range("M22").Value = 27.26348539453592
$ws.Range("N22").Value = 18.94349669364409
$ws.Range("O22").Value = 25.81550220148414
$ws.Range("C23").Value = 8.888944996140378
$ws.Range("D23").Value = 4.777167117680089
$ws.Range("E23").Value = 11.1787010510516
$ws.Range("F23").Value = 29.84416413475932
$ws.Range("G23").Value = 3.61571795627452
$ws.Range("I23").Value = 25.84446652150768
$ws.Range("L23").Value = 8.77991664990367
$ws.Range("M23").Value = 26.96937776498636
$ws.Range("N23").Value = 18.81508842305514
$ws.Range("O23").Value = 25.72872151920502
$ws.Range("C24").Value = 8.951593390345449
$ws.Range("D24").Value = 4.800753570666223
$ws.Range("E24").Value = 11.28729436977545
$ws.Range("F24").Value = 29.28510408562004
$ws.Range("G24").Value = 3.620460145967725
$ws.Range("I24").Value = 25.60978357299923
$ws.Range("L24").Value = 8.845610999472575
$ws.Range("M24").Value = 25.82473756849567
$ws.Range("N24").Value = 18.32327912899035
$ws.Range("O24").Value = 25.41099214046183
$ws.Range("C25").Value = 9.025234056847825
$ws.Range("D25").Value = 4.827543176061017
$ws.Range("E25").Value = 11.41215393148494
$ws.Range("F25").Value = 28.70583468650284
$ws.Range("G25").Value = 3.625954771963517
$ws.Range("I25").Value = 25.38348382966527
$ws.Range("L25").Value = 8.921043200693182
$ws.Range("M25").Value = 24.52995853255413
$ws.Range("N25").Value = 17.38162619306493
$ws.Range("O25").Value = 25.0930518843957
